$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.961.65'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.959.18'
$ws.Range('E3').Value = '  +2.50%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '353.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.563'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.633'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.75'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0899'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.26%  '
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.91'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.76%  '
$ws.Range('D15').Value = '3.426.19'
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('D16').Value = '2.969.15'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '52.095.70'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.22%  '
$ws.Range('D22').Value = '0.0₃0994'
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '272.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('E26').Value = '  +10.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('E29').Value = '  +18.27%  '
$ws.Range('E30').Value = '  +22.16%  '
$ws.Range('E31').Value = '  +1.29%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.39'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.78%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '37.86'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '53.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0451'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.42'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('B38').Value = 'Toncoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.88'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -17.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.94'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.06'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.68%  '
$ws.Range('E42').Value = '  +2.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.99%  '
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('D47').Value = '2.172.19'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '114.10'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.95%  '
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0341'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.936'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.35%  '
